$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 195; this shifts the existing rows 195..302
# down to 196..303, preserving all of their original data untouched.
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with the new data point.
$ws.Cells.Item(195, 1).Value2  = 4
$ws.Cells.Item(195, 2).Value2  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(195, 3).Value2  = 'Los Lagos'
$ws.Cells.Item(195, 4).Value2  = 44719
$ws.Cells.Item(195, 5).Value2  = 10
$ws.Cells.Item(195, 6).Value2  = 100114014
$ws.Cells.Item(195, 7).Value2  = 'Betarraga'
$ws.Cells.Item(195, 8).Value2  = 'Sin especificar'
$ws.Cells.Item(195, 9).Value2  = 'Primera'
$ws.Cells.Item(195, 10).Value2 = 900
$ws.Cells.Item(195, 11).Value2 = 1200
$ws.Cells.Item(195, 12).Value2 = 1200
$ws.Cells.Item(195, 13).Value2 = 1200
$ws.Cells.Item(195, 14).Value2 = '$/paquete 5 unidades'
$ws.Cells.Item(195, 15).Value2 = 'Región del Maule'
$ws.Cells.Item(195, 16).Value2 = 240
$ws.Cells.Item(195, 17).Value2 = 5
$ws.Cells.Item(195, 18).Value2 = 'Hortaliza'
